$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the two new number-format styles in the same order they appear
#     in the target styles.xml (cellXfs index 1 = date/14, index 2 = text/49)
#     by formatting one cell and copy/paste-special'ing the format onto the
#     rest of the range so every cell shares a single style record instead
#     of each getting its own duplicate xf entry. ---

# cellXfs index 1: numFmtId 14 (built-in date), applied to the (currently
# empty) helper cells C8:C10.
$ws.Range("C8").NumberFormat = "mm-dd-yy"
$ws.Range("C8").Copy()
$ws.Range("C9:C10").PasteSpecial(-4122)

# cellXfs index 2: numFmtId 49 (Text), applied to the date column of the
# data rows C2:C4 so the date strings stay text instead of being coerced
# into date serials.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Copy()
$ws.Range("C3:C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Range("A2").Value2 = "I love the new car policy, it was much cheaper than my previous plan"
$ws.Range("B2").Value2 = "Steve Doe"
$ws.Range("C2").Value2 = "9/10/2021"
$ws.Range("D2").Value2 = 45
$ws.Range("E2").Value2 = "Car"
$ws.Range("F2").Value2 = "Vehicle"

# --- Row 3 ---
$ws.Range("A3").Value2 = "Great service, Jake was really fast and helpful and helping me with my claim"
$ws.Range("B3").Value2 = "John Doe"
$ws.Range("C3").Value2 = "3/2/2016"
$ws.Range("D3").Value2 = 100
$ws.Range("E3").Value2 = "Whole"
$ws.Range("F3").Value2 = "Life"

# --- Row 4 (new row) ---
$ws.Range("A4").Value2 = "I had a really unpleasant time with the new app, it was bad"
$ws.Range("B4").Value2 = "Sarah Doe"
$ws.Range("C4").Value2 = "6/25/2011"
$ws.Range("D4").Value2 = 250
$ws.Range("E4").Value2 = "Homeowner"
$ws.Range("F4").Value2 = "Home"

# --- Column widths ---
$ws.Columns.Item(3).ColumnWidth = 25.140625

# --- View / selection ---
$ws.Range("C5").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1
